$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 5: A5 = URL (hyperlink), B5 = Title (large bold font)
$url = "https://mccormickml.com/2019/05/14/BERT-word-embeddings-tutorial/"
$title = "BERT Word Embeddings Tutorial"

$ws.Range("A5").Value = $url
$ws.Range("B5").Value = $title

$ws.Hyperlinks.Add($ws.Range("A5"), $url, "", "", $url)

# Style A5 like the other link cells in column A (blue, wrap text) but in Arial
$a5 = $ws.Range("A5")
$a5.Font.Name = "Arial"
$a5.Font.Size = 10
$a5.Font.Color = 255
$a5.WrapText = $true

# Style B5 as a large bold Arial title
$b5 = $ws.Range("B5")
$b5.Font.Name = "Arial"
$b5.Font.Size = 22
$b5.Font.Bold = $true
$b5.WrapText = $true

# Update row height for the new row to fit the larger font
$ws.Rows.Item(5).RowHeight = 26.8

# Move the active selection to A6, matching the post-edit cursor position
$ws.Range("A6").Select()
